$d = $word.ActiveDocument

# --- 1) Merge "Arduino External Power Supply" + " x2" into one run ---
$d.Content.Find.Execute("Arduino External Power Supply x2", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Arduino External Power Supply x2", 2) | Out-Null

# --- 2) Merge "(" + "w/ " + "9V Battery Adapter)" into one run ---
$d.Content.Find.Execute("(w/ 9V Battery Adapter)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "(w/ 9V Battery Adapter)", 2) | Out-Null

# --- 3) Merge the "especially in our project is that" run sequence (drops the
#        proofErr gramStart/gramEnd bookmarks around "project" as a side effect
#        of Find/Replace recombining the runs it touches) ---
$d.Content.Find.Execute("especially in our project is that", $false, $false, $false, $false, $false, `
    $true, 1, $false, "especially in our project is that", 2) | Out-Null

# --- 4) Replace the two trailing empty paragraphs with a "GitHub Repository
#        Link" heading paragraph + the URL paragraph, leaving the final
#        trailing empty paragraph (before the sectPr) untouched. ---

# Find the run of trailing empty paragraphs (right before the very last one,
# which must stay empty) by scanning back from the end of the document.
$total = $d.Paragraphs.Count
$lastEmptyIdx = 0
$i = $total
while ($i -ge 1) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "") {
        $lastEmptyIdx = $i
    } else {
        break
    }
    $i = $i - 1
}

# $lastEmptyIdx .. $total are the trailing empty paragraphs. Keep the very
# last one empty; use the two immediately before it for the new content.
$targetIdx = $total - 2

# Collapse the pair of empty paragraphs we'll replace down to a single
# paragraph by deleting the one right after it.
$d.Paragraphs($targetIdx + 1).Range.Delete() | Out-Null

$p = $d.Paragraphs($targetIdx)
$rng = $p.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>GitHub Repository Link</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>https://github.com/surenb1/Heat-Seeking-Robot-Arduino</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml) | Out-Null

Write-Output "done"
